$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-22 18:31:00"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
